$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110, column C ("Numéro de page") no longer has a value ("NA") -
# clear it back to a blank cell (keeping the cell reference, no style residue).
$ws.Range("C110").NumberFormat = "@"
$ws.Range("C110").Value = ""
$ws.Range("C110").Style = "Normal"

# New row 111: the entry that used to carry "NA" in C110 now lives on its own
# row for 2025-05-15.
$ws.Range("A111").NumberFormat = "@"
$ws.Range("A111").Value = "2025-05-15"
$ws.Range("A111").Style = "Normal"

$ws.Range("B111").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C111").Value = "NA"
$ws.Range("D111").Value = 1
